# Update the "dSF" (column F) values on the active sheet to reflect the
# repulled/recomputed data, per the commit: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 9
    3  = 1
    4  = -2
    5  = -2
    8  = -5
    9  = 2
    10 = 1
    11 = -2
    12 = 5
    13 = 12
    14 = 4
    15 = 5
    16 = 6
    17 = -4
    18 = 10
    19 = 3
    20 = 2
    21 = -1
    22 = -3
    23 = 7
    25 = 6
    26 = 1
    27 = 1
    29 = -3
    30 = 6
    31 = -2
    33 = -1
    34 = 2
    35 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
